$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Interpolation_instructions")
$ws2 = $wb.Worksheets.Item("Trend_instructions")

# Add new shared string header on sheet2 column K1
$ws2.Range("K1").Value = "override_normalization"

# Update F5 value from 1 to 2 on sheet2
$ws2.Range("F5").Value = 2

# Auto-fit the new column so its width reflects the new header text
$ws2.Columns.Item(10).EntireColumn.AutoFit()

# Update selections (view state) to match diff; leave Trend_instructions
# (sheet2) as the active/selected sheet, matching the saved workbook state
$ws1.Activate()
$ws1.Range("I1").Select()
$ws2.Activate()
$ws2.Range("F14").Select()

